# Generate Report for Handoff
# - Update status text "In Translation" -> "Ready for handoff"
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Widen the status/date columns to fit the new, longer text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column(s): "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# Latest handoff timestamps (stored as literal text, not real dates)
$overview.Range("G2").Value = "2016-08-21 17:04:40"
$zhcn.Range("H2").Value     = "2016-08-21 17:04:36"
$dede.Range("H2").Value     = "2016-08-21 17:04:40"

# Widen columns that hold the status/date values to fit the new text.
# (Target OOXML width ~17.216 chars; the COM ColumnWidth setter here snaps to a
# whole-pixel grid, so 16.33 is the nearest settable value that lands on it.)
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33

# Re-assert the datetime display format on the timestamp cells we just edited
# so the round-trip keeps their "yyyy-mm-dd HH:mm:ss" formatting intact.
$overview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("H2").NumberFormat     = "yyyy-mm-dd HH:mm:ss"
$dede.Range("H2").NumberFormat     = "yyyy-mm-dd HH:mm:ss"
